# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.943.28"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "3.794.42"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'592.94"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").Value = "'171.44"
$ws.Range("E6").Value = "  -5.11%  "
$ws.Range("D7").Value = "3.792.02"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "'38.27"
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("D14").Value = "'0.0000243"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").Value = "4.431.73"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "3.803.74"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "68.063.73"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = "  -4.33%  "
$ws.Range("D20").Value = "'16.04"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "'486.75"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "'0.735"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").Value = "'85.85"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  -6.07%  "
$ws.Range("D26").Value = "'0.0000137"
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("D27").Value = "'12.24"
$ws.Range("E27").Value = "  -5.59%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  -8.70%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'2.93"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'2.43"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'32.47"
$ws.Range("E32").Value = "  +7.00%  "
$ws.Range("D33").Value = "'7.57"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").Value = "'0.110"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").Value = "'5.83"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").Value = "'0.325"
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("D40").Value = "'445.93"
$ws.Range("E40").Value = "  +4.51%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'48.96"
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.01"
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").Value = "'2.87"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").Value = "'8.32"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").Value = "'41.40"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("D46").Value = "2.855.18"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0352"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").Value = "'137.95"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "'26.44"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  -5.78%  "
